$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width changes ---
$ws.Range("B1").ColumnWidth = 25.3775510204082
$ws.Range("BX1").ColumnWidth = 5.12755102040816
$ws.Range("BY1").ColumnWidth = 5.66836734693878
$ws.Range("BZ1").ColumnWidth = 26.7295918367347
$ws.Range("CA1").ColumnWidth = 26.3214285714286
$ws.Range("CB1:AMK1").EntireColumn.ColumnWidth = 8.23469387755102

# --- New attendance cell values (lecture 3: 4 new date columns BX..CA) ---
$ws.Range("BZ2").Value = "P"
$ws.Range("CA2").Value = "Leave"
$ws.Range("BX3").Value = "Leave"
$ws.Range("BY3").Value = "Leave"
$ws.Range("BZ3").Value = "Leave"
$ws.Range("CA3").Value = "Leave"
$ws.Range("BY4").Value = "P"
$ws.Range("BZ4").Value = "P"
$ws.Range("CA4").Value = "P"
$ws.Range("BX5").Value = "Leave"
$ws.Range("BY5").Value = "P"
$ws.Range("BZ5").Value = "Leave"
$ws.Range("CA5").Value = "P"
$ws.Range("BY6").Value = "Leave"
$ws.Range("BZ6").Value = "Leave"
$ws.Range("CA6").Value = "Leave"
$ws.Range("BY7").Value = "Leave"
$ws.Range("BZ7").Value = "Leave"
$ws.Range("CA7").Value = "P"
$ws.Range("BY8").Value = "P"
$ws.Range("BZ8").Value = "Leave"
$ws.Range("CA8").Value = "Leave"
$ws.Range("BZ9").Value = "P"
$ws.Range("BY11").Value = "P"
$ws.Range("BZ11").Value = "P"
$ws.Range("CA11").Value = "P"
$ws.Range("BY12").Value = "P"
$ws.Range("BZ12").Value = "P"
$ws.Range("CA12").Value = "P"
$ws.Range("BY13").Value = "P"
$ws.Range("BZ13").Value = "Leave"
$ws.Range("CA13").Value = "Leave"
$ws.Range("BY14").Value = "P"
$ws.Range("BZ14").Value = "P"
$ws.Range("BY15").Value = "Leave"
$ws.Range("BZ15").Value = "Leave"
$ws.Range("CA15").Value = "P"
$ws.Range("BY16").Value = "P"
$ws.Range("BZ16").Value = "P"
$ws.Range("CA16").Value = "P"
$ws.Range("BZ17").Value = "Leave"
$ws.Range("BY18").Value = "P"
$ws.Range("BZ18").Value = "P"
$ws.Range("CA18").Value = "Leave"
$ws.Range("BX20").Value = "Leave"
$ws.Range("BY20").Value = "Leave"
$ws.Range("BZ20").Value = "Leave"
$ws.Range("CA20").Value = "P"
$ws.Range("BY21").Value = "Leave"
$ws.Range("BZ21").Value = "P"
$ws.Range("CA21").Value = "Leave"
$ws.Range("BY23").Value = "P"
$ws.Range("BZ23").Value = "P"
$ws.Range("CA23").Value = "P"
$ws.Range("BY24").Value = "P"
$ws.Range("BZ24").Value = "P"
$ws.Range("CA24").Value = "Leave"
$ws.Range("BY25").Value = "P"
$ws.Range("BZ25").Value = "P"
$ws.Range("CA25").Value = "Leave"
$ws.Range("BY26").Value = "P"
$ws.Range("BZ26").Value = "P"
$ws.Range("CA26").Value = "P"
$ws.Range("BY27").Value = "Leave"
$ws.Range("BZ27").Value = "P"
$ws.Range("CA27").Value = "Leave"
$ws.Range("BY28").Value = "P"
$ws.Range("BZ28").Value = "P"
$ws.Range("CA28").Value = "Leave"
$ws.Range("BY29").Value = "P"
$ws.Range("BZ29").Value = "P"
$ws.Range("BY30").Value = "Leave"
$ws.Range("BZ30").Value = "P"
$ws.Range("CA30").Value = "Leave"
$ws.Range("BY31").Value = "Leave"
$ws.Range("BZ31").Value = "Leave"
$ws.Range("CA31").Value = "P"
$ws.Range("BY32").Value = "P"
$ws.Range("BZ32").Value = "P"
$ws.Range("CA32").Value = "P"
$ws.Range("BY33").Value = "Leave"
$ws.Range("BZ33").Value = "P"
$ws.Range("CA33").Value = "Leave"
$ws.Range("BZ35").Value = "Leave"
$ws.Range("BY36").Value = "P"
$ws.Range("BZ36").Value = "P"
$ws.Range("CA36").Value = "P"
$ws.Range("BY37").Value = "P"
$ws.Range("BZ37").Value = "P"
$ws.Range("CA37").Value = "P"
$ws.Range("BY38").Value = "P"
$ws.Range("BZ38").Value = "Leave"
$ws.Range("BY39").Value = "P"
$ws.Range("BZ39").Value = "P"
$ws.Range("CA39").Value = "Leave"
$ws.Range("BY41").Value = "Leave"
$ws.Range("BZ41").Value = "P"
$ws.Range("CA41").Value = "P"
$ws.Range("BY43").Value = "P"
$ws.Range("BZ43").Value = "P"
$ws.Range("CA43").Value = "P"
$ws.Range("BZ44").Value = "P"
$ws.Range("BY45").Value = "P"
$ws.Range("BZ45").Value = "P"
$ws.Range("CA45").Value = "P"
$ws.Range("BZ46").Value = "P"
$ws.Range("BY47").Value = "P"
$ws.Range("BZ47").Value = "P"
$ws.Range("CA47").Value = "P"
$ws.Range("BZ48").Value = "Leave"
$ws.Range("CA48").Value = "Leave"
$ws.Range("BY50").Value = "P"
$ws.Range("BZ50").Value = "P"
$ws.Range("CA50").Value = "P"
$ws.Range("BZ51").Value = "Leave"

# --- Selection / view state ---
$null = $ws.Range("CA36").Select()
